# Update countries & provincias Spain
# - Swap the "Etiopia"/"Portugal" rows (row 51 becomes Etiopia, row 52 becomes Portugal)
# - Refresh the "Datos actualizados..." timestamp
# - Update the daily COVID figures (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes hoy, Muertes) for the affected countries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header banner: refresh the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 20:37"

# --- Country name swap: Etiopia now ranks above Portugal ---
$ws.Range("A51").Value = "Etiopia"
$ws.Range("A52").Value = "Portugal"

# --- Row 4 (Estados Unidos) ---
$ws.Range("B4").Value = 6764598
$ws.Range("C4").Value = 15309
$ws.Range("D4").Value = 4043637
$ws.Range("E4").Value = 2521339
$ws.Range("G4").Value = 622
$ws.Range("H4").Value = 199622

# --- Row 5 (India) ---
$ws.Range("B5").Value = 5017930
$ws.Range("C5").Value = 91016
$ws.Range("D5").Value = 3939048
$ws.Range("E5").Value = 996794
$ws.Range("G5").Value = 1280
$ws.Range("H5").Value = 82088

# --- Row 16 (Francia) ---
$ws.Range("B16").Value = 395104
$ws.Range("C16").Value = 7852
$ws.Range("D16").Value = 89891
$ws.Range("E16").Value = 274214
$ws.Range("G16").Value = 49
$ws.Range("H16").Value = 30999

# --- Row 25 (Alemania) ---
$ws.Range("B25").Value = 264375
$ws.Range("C25").Value = 1154
$ws.Range("E25").Value = 17383

# --- Row 41 (Marruecos) ---
$ws.Range("B41").Value = 90324
$ws.Range("C41").Value = 2121
$ws.Range("D41").Value = 71047
$ws.Range("E41").Value = 17629
$ws.Range("G41").Value = 34
$ws.Range("H41").Value = 1648

# --- Row 51 (now Etiopia) ---
$ws.Range("B51").Value = 65486
$ws.Range("C51").Value = 700
$ws.Range("D51").Value = 25988
$ws.Range("E51").Value = 38463
$ws.Range("G51").Value = 13
$ws.Range("H51").Value = 1035

# --- Row 52 (now Portugal) ---
$ws.Range("B52").Value = 65021
$ws.Range("C52").Value = 425
$ws.Range("D52").Value = 44362
$ws.Range("E52").Value = 18784
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 1875

# --- Row 73 (Irlanda) ---
$ws.Range("B73").Value = 31549
$ws.Range("C73").Value = 357
$ws.Range("E73").Value = 6398
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 1787

# --- Row 77 (Libano) ---
$ws.Range("D77").Value = 9216
$ws.Range("E77").Value = 15981

# --- Row 99 (Guayana Francesa) ---
$ws.Range("B99").Value = 9578
$ws.Range("C99").Value = 26
$ws.Range("D99").Value = 9202
$ws.Range("E99").Value = 312
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 64

# --- Row 141 (Mali) ---
$ws.Range("B141").Value = 2940
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 2302
$ws.Range("E141").Value = 510

# --- Row 142 (Reunion) ---
$ws.Range("B142").Value = 2902
$ws.Range("C142").Value = 30
$ws.Range("E142").Value = 1574

# --- Row 144 (Sudan del Sur) ---
$ws.Range("B144").Value = 2592
$ws.Range("C144").Value = 5
$ws.Range("E144").Value = 1253

# --- Row 152 (Yemen) ---
$ws.Range("B152").Value = 2016
$ws.Range("C152").Value = 3
$ws.Range("D152").Value = 1219
$ws.Range("E152").Value = 214

# --- Row 180 (Eritrea) ---
$ws.Range("B180").Value = 364
$ws.Range("C180").Value = 3
$ws.Range("E180").Value = 60

# --- Row 192 (Curazao) ---
$ws.Range("B192").Value = 169
$ws.Range("C192").Value = 8
$ws.Range("D192").Value = 61
$ws.Range("E192").Value = 107
